# Add daily power records
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comforter-cda")

# New data rows to append: Date, Start Time, End Time
$newRows = @(
    @{ Row = 117; Date = 43441; Start = $null; End = $null },
    @{ Row = 118; Date = 43442; Start = 0.43263888888888885; End = 0.45208333333333334 },
    @{ Row = 119; Date = 43442; Start = 0.59861111111111109; End = 0.6743055555555556 },
    @{ Row = 120; Date = $null; Start = $null; End = $null },
    @{ Row = 121; Date = $null; Start = $null; End = $null }
)

foreach ($r in $newRows) {
    $row = $r.Row
    if ($null -ne $r.Date) {
        $ws.Cells.Item($row, 1).Value = $r.Date
    }
    if ($null -ne $r.Start) {
        $ws.Cells.Item($row, 2).Value = $r.Start
    }
    if ($null -ne $r.End) {
        $ws.Cells.Item($row, 3).Value = $r.End
    }
    $ws.Cells.Item($row, 4).Formula = "=(C$row-B$row)* 1440"
    $ws.Cells.Item($row, 5).Formula = "=IF(C$row>B$row, (C$row-B$row)*1440, (B$row-C$row)*1440)"
    $ws.Cells.Item($row, 6).Formula = "=ABS((C$row-B$row)*1440)"
}

# Resize the table to include the new rows
$table = $ws.ListObjects.Item("comforter_cda_table")
$table.Resize($ws.Range("A1:F121"))

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 105
$ws.Range("C120").Select()
